$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 0.0165316589608312
$ws.Range("C3").Value = 0.1596747772390028
$ws.Range("C4").Value = 0.8924876052044113
$ws.Range("C5").Value = 0.7209875032280612
$ws.Range("C6").Value = 0.9952126454880028
$ws.Range("C7").Value = -0.001024306990703626
$ws.Range("C8").Value = -0.09342248582433159
$ws.Range("C9").Value = -0.2773457422942077
$ws.Range("C10").Value = 0.673185173740518
$ws.Range("C11").Value = -0.8972713613521036
$ws.Range("C12").Value = 0.7655094493981847
$ws.Range("C13").Value = 0.1599476620886939
$ws.Range("C14").Value = 0.3444026411998959
$ws.Range("C15").Value = 0.07589299699726311
$ws.Range("C16").Value = 2.453088300160251
$ws.Range("C17").Value = -0.1549791377016192
$ws.Range("C18").Value = 0.7546548098942604
$ws.Range("C19").Value = 0.2903360318096634
$ws.Range("C20").Value = 0.6557382247212946
$ws.Range("C21").Value = 0.8294043729654725
$ws.Range("C22").Value = 0.7192499639843549
$ws.Range("C23").Value = 0.4887190027762116
$ws.Range("C24").Value = 1.715815238194143
$ws.Range("C25").Value = 1.264219810813607
$ws.Range("C26").Value = 0.9808549446337534
$ws.Range("C27").Value = 0.3866366398948996
$ws.Range("C28").Value = 0.4914602873909288
$ws.Range("C29").Value = 0.8003563227094359
$ws.Range("C30").Value = 0.7361494804832538
$ws.Range("C31").Value = 0.7587805559099209
